$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so numeric-looking strings
# (e.g. "1.004", "0.000008032") are stored as text, matching the
# original inlineStr cell type instead of being parsed as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.194.53"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.670.97"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "210.70"
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("D6").Value = "0.5214"
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "0.2626"
$ws.Range("D9").Value = "0.06334"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "21.25"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "0.07556"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "1.675.10"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "4.446"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").Value = "0.5502"
$ws.Range("E14").Value = "  -4.11%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "0.000008032"
$ws.Range("E15").Value = "  -2.20%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "66.53"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Value = "26.204.67"
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "4.754"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").Value = "187.11"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("E21").Value = "  -4.13%  "
$ws.Range("D22").Value = "6.220"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "150.09"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("D25").Value = "0.1245"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").Value = "7.504"
$ws.Range("E26").Value = "  -3.37%  "
$ws.Range("D27").Value = "15.85"
$ws.Range("E27").Value = "  +0.85%  "
$ws.Range("D28").Value = "0.06333"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").Value = "1.359"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "1.280"
$ws.Range("E30").Value = "  -2.49%  "
$ws.Range("D31").Value = "3.526"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "3.417"
$ws.Range("E32").Value = "  -4.15%  "
$ws.Range("D33").Value = "1.647"
$ws.Range("E33").Value = "  -1.98%  "
$ws.Range("E34").Value = "  -1.86%  "
$ws.Range("D35").Value = "0.6041"
$ws.Range("E35").Value = "  -1.04%  "
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.112.52"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "6.141"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("D41").Value = "0.8677"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").Value = "1.005"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").Value = "100.43"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").Value = "1.824.66"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").Value = "0.00000000107"
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("D46").Value = "55.65"
$ws.Range("E46").Value = "  -2.89%  "
$ws.Range("D47").Value = "1.002"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "8.095"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").Value = "0.05240"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").Value = "0.4245"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").Value = "5.939"
$ws.Range("E51").Value = "  -1.00%  "

# Remove the temporary text-number-format styling so the cells keep
# their original (default/no explicit style) appearance.
$ws.Range("D2:D51").ClearFormats()
